# Integration plan reworked bottom-up using the dependency tree:
# the old rows 5-7 are removed, and the remaining rows (2-4) are
# updated to reflect the new step dependencies.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Remove the now-obsolete rows (old steps 4, 5, 6) - this also shrinks
# the table/dimension/used-range from A1:H7 down to A1:H4.
$ws.Rows("5:7").Delete()

# Row 2 (Step 1 - TransponderReceiver): add dependency marker in C,
# and flip the remaining downstream markers from "S" to "X".
$ws.Range("C2").Value = "S"
$ws.Range("F2").Value = "X"
$ws.Range("G2").Value = "X"
$ws.Range("H2").Value = "X"

# Row 3 (Step 2 - Decoder): add a "T" dependency marker in C,
# change D from "T" to "X", and flip remaining markers to "X".
$ws.Range("C3").Value = "T"
$ws.Range("D3").Value = "X"
$ws.Range("G3").Value = "X"
$ws.Range("H3").Value = "X"

# Row 4 (Step 3 - TrackHandler): add "T" in B and "X" in C,
# change D from "T" to "X", and flip H from "S" to "X".
$ws.Range("B4").Value = "T"
$ws.Range("C4").Value = "X"
$ws.Range("D4").Value = "X"
$ws.Range("H4").Value = "X"

# Restore the view: selection moves to G9, and the previous
# topLeftCell scroll-freeze anchor is cleared by reselecting.
$ws.Range("G9").Select()
